# Generate Report for Archive
#
# 1) Update the "Ready for handoff" status text to "In Translation" everywhere
#    it appears (Overview!E2:F3 and the Status column (C2:C3) on each of the
#    per-language sheets).
# 2) Shrink the width of the columns that hold that status text (Overview
#    columns E/F, and column C on each language sheet) to match the
#    narrower report layout.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# Narrow the zh-cn / de-de columns on the Overview sheet.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = "In Translation"
    $ws.Range("C3").Value = "In Translation"

    # Narrow the Status column to match.
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
